$p = $ppt.ActivePresentation

# --- Slide 2: Implementation ---
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Implementation"
$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "PHP"
[void]$body2.InsertAfter("`rMySQL")
[void]$body2.InsertAfter("`rHTML")
[void]$body2.InsertAfter("`rCSS")
[void]$body2.InsertAfter("`rJavaScript")
# trailing blank bullet line (no text) after the last item
[void]$body2.InsertAfter("`rx")
$body2.Characters($body2.Length, 1).Text = ""

# --- Slide 3: Screenshots ---
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Screenshots"

# --- Slide 4: Issues ---
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Issues"
$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "Make the editor dynamic"
[void]$body4.InsertAfter("`rHow to represent timeline information")
[void]$body4.InsertAfter("`rHow to distribute application")
# trailing blank line after the last item
[void]$body4.InsertAfter("`rx")
$body4.Characters($body4.Length, 1).Text = ""

# --- Slide 5: Conclusions ---
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusions"
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Learned how to develop full scale web application"
[void]$body5.InsertAfter("`rLearned how to represent data with JSON")
[void]$body5.InsertAfter("`rLearned about deploying ")
[void]$body5.InsertAfter("web applications")
